$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks before clearing so no stale relationships remain
$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# Header row
$ws.Range('A1').Value = '取得日時'
$ws.Range('B1').Value = 'タイトル'
$ws.Range('C1').Value = 'カテゴリ'
$ws.Range('D1').Value = '価格'
$ws.Range('E1').Value = '締切'
$ws.Range('F1').Value = 'URL'
$ws.Range('G1').Value = '優先度スコア'
$ws.Range('H1').Value = 'スキル概要'

# Row 2
$ws.Range('A2').Value = '2025-09-11 12:33:56'
$ws.Range('B2').Value = '【業務委託/副業可】AI SaaS開発を牽引するCTO候補を募集'
$ws.Range('C2').Value = 'システム開発'
$ws.Range('D2').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E2').Value = '期限情報なし'
$ws.Range('F2').Value = 'https://www.lancers.jp/work/detail/5391756'
$ws.Range('G2').Value = 375
$ws.Range('H2').Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Range('A3').Value = '2025-09-11 12:33:56'
$ws.Range('B3').Value = '【AIで開発生産性を革新】AI活用推進エンジニア募集(副業・業務委託)'
$ws.Range('C3').Value = 'システム開発'
$ws.Range('D3').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E3').Value = '期限情報なし'
$ws.Range('F3').Value = 'https://www.lancers.jp/work/detail/5391761'
$ws.Range('G3').Value = 375
$ws.Range('H3').Value = '🔥AI,Ai ◆開発'

# Row 4
$ws.Range('A4').Value = '2025-09-11 12:33:56'
$ws.Range('B4').Value = '競馬AIの開発ができる方、もしくはすでに開発済みの方'
$ws.Range('C4').Value = 'システム開発'
$ws.Range('D4').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E4').Value = '期限情報なし'
$ws.Range('F4').Value = 'https://www.lancers.jp/work/detail/5391744'
$ws.Range('G4').Value = 375
$ws.Range('H4').Value = '🔥AI,Ai ◆開発'

# Row 5
$ws.Range('A5').Value = '2025-09-11 12:33:56'
$ws.Range('B5').Value = '【AI技術顧問/戦略アドバイザー募集】最先端AIで事業の非連続な成長を牽引するエキスパート求む'
$ws.Range('C5').Value = 'システム開発'
$ws.Range('D5').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E5').Value = '期限情報なし'
$ws.Range('F5').Value = 'https://www.lancers.jp/work/detail/5391776'
$ws.Range('G5').Value = 310
$ws.Range('H5').Value = '🔥AI,Ai'

# Row 6
$ws.Range('A6').Value = '2025-09-11 12:33:56'
$ws.Range('B6').Value = '【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)'
$ws.Range('C6').Value = 'システム開発'
$ws.Range('D6').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E6').Value = '期限情報なし'
$ws.Range('F6').Value = 'https://www.lancers.jp/work/detail/5391607'
$ws.Range('G6').Value = 155
$ws.Range('H6').Value = '◆開発,Node.js'

# Row 7
$ws.Range('A7').Value = '2025-09-11 12:33:56'
$ws.Range('B7').Value = '【急募】SharePoint+Power Platformでの不動産賃貸管理システム構築'
$ws.Range('C7').Value = 'システム開発'
$ws.Range('D7').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E7').Value = '期限情報なし'
$ws.Range('F7').Value = 'https://www.lancers.jp/work/detail/5391490'
$ws.Range('G7').Value = 60
$ws.Range('H7').Value = '◇管理'

# Row 8
$ws.Range('A8').Value = '2025-09-11 12:33:56'
$ws.Range('B8').Value = '【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!'
$ws.Range('C8').Value = 'システム開発'
$ws.Range('D8').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E8').Value = '期限情報なし'
$ws.Range('F8').Value = 'https://www.lancers.jp/work/detail/5371747'
$ws.Range('G8').Value = 48
$ws.Range('H8').Value = '◆コンサル'

# Row 9
$ws.Range('A9').Value = '2025-09-11 12:33:56'
$ws.Range('B9').Value = '【急募・再掲載】自社アプリのデバッグ・バグチェック業務依頼 ※NDA締結必須'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '~ 5,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('F9').Value = 'https://www.lancers.jp/work/detail/5391844'
$ws.Range('G9').Value = 30
$ws.Range('H9').Value = '◇アプリ'

# Row 10
$ws.Range('A10').Value = '2025-09-11 12:33:56'
$ws.Range('B10').Value = '初回 【フルリモート】フリーランスエンジニア募集'
$ws.Range('C10').Value = 'システム開発'
$ws.Range('D10').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E10').Value = '期限情報なし'
$ws.Range('F10').Value = 'https://www.lancers.jp/work/detail/5391489'
$ws.Range('G10').Value = 25

# Row 11
$ws.Range('A11').Value = '2025-09-11 12:33:56'
$ws.Range('B11').Value = '要件定義や基本設計ができる方(1人月、約2年アサイン予定)'
$ws.Range('C11').Value = 'システム開発'
$ws.Range('D11').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E11').Value = '期限情報なし'
$ws.Range('F11').Value = 'https://www.lancers.jp/work/detail/5391221'
$ws.Range('G11').Value = 25

# Row 12
$ws.Range('A12').Value = '2025-09-11 12:33:56'
$ws.Range('B12').Value = '【講師募集】Gensparkを使ったWEB構築チュートリアル募集'
$ws.Range('C12').Value = 'システム開発'
$ws.Range('D12').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E12').Value = '期限情報なし'
$ws.Range('F12').Value = 'https://www.lancers.jp/work/detail/5390165'
$ws.Range('G12').Value = 13

# Row 13
$ws.Range('A13').Value = '2025-09-11 12:33:56'
$ws.Range('B13').Value = '【急募】Googleアナリティクス連携の専門家を探しています'
$ws.Range('C13').Value = 'システム開発'
$ws.Range('D13').Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range('E13').Value = '期限情報なし'
$ws.Range('F13').Value = 'https://www.lancers.jp/work/detail/5391267'
$ws.Range('G13').Value = 10

# Recreate hyperlinks for F2:F13 in order, so relationship ids map rId1->F2 ... rId12->F13
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5391756') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5391761') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5391744') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5391776') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5391607') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5391490') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5371747') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5391844') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F10'), 'https://www.lancers.jp/work/detail/5391489') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F11'), 'https://www.lancers.jp/work/detail/5391221') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F12'), 'https://www.lancers.jp/work/detail/5390165') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F13'), 'https://www.lancers.jp/work/detail/5391267') | Out-Null

# Restore the Hyperlink cell style (Add() leaves a duplicated style index otherwise)
$ws.Range('F2').Style = 'Hyperlink'
$ws.Range('F3').Style = 'Hyperlink'
$ws.Range('F4').Style = 'Hyperlink'
$ws.Range('F5').Style = 'Hyperlink'
$ws.Range('F6').Style = 'Hyperlink'
$ws.Range('F7').Style = 'Hyperlink'
$ws.Range('F8').Style = 'Hyperlink'
$ws.Range('F9').Style = 'Hyperlink'
$ws.Range('F10').Style = 'Hyperlink'
$ws.Range('F11').Style = 'Hyperlink'
$ws.Range('F12').Style = 'Hyperlink'
$ws.Range('F13').Style = 'Hyperlink'

# Column B width -> 49 characters (ColumnWidth input empirically mapped to hit an exact integer)
$ws.Columns.Item(2).ColumnWidth = 48.14

Write-Output $ws.UsedRange.Address()